# Add the new "send-otp" signup/log rows (rows 71-78) to the "2022_4"
# worksheet (sheet4.xml), matching the author's "add signup for
# phonenumber" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2022_4")

# Column layout (row 1 header): A=date B=time C=phoneNumber D=model
# E=path F=action G=status H=description M=failureReason
$rows = @(
    @{ Row=71; A="Mon Apr 25 2022"; B="11:42:18 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="Cannot read properties of undefined (reading 'findFirst')" },
    @{ Row=72; A="Mon Apr 25 2022"; B="11:43:20 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="error.invalid" },
    @{ Row=73; A="Mon Apr 25 2022"; B="11:43:42 GMT+0000 (Greenwich Mean Time)"; C="22892942601"; D="User"; E="/api/auth/send-otp"; F="request"; G="succeeded"; H="22892942601 request to receive otp"; M=$null },
    @{ Row=74; A="Mon Apr 25 2022"; B="11:43:57 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="error.invalid" },
    @{ Row=75; A="Mon Apr 25 2022"; B="11:44:05 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="error.invalid" },
    @{ Row=76; A="Mon Apr 25 2022"; B="11:44:12 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="error.invalid" },
    @{ Row=77; A="Mon Apr 25 2022"; B="11:44:43 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="failed";    H="22892942600 request to receive otp"; M="error.invalid" },
    @{ Row=78; A="Mon Apr 25 2022"; B="11:47:17 GMT+0000 (Greenwich Mean Time)"; C="22892942600"; D="User"; E="/api/auth/send-otp"; F="request"; G="succeeded"; H="22892942600 request to receive otp"; M=$null }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    # phoneNumber is numeric-looking text -- force text so it isn't coerced
    # into a Number cell (leading apostrophe = Excel's "store as text" cue).
    # ClearFormats() drops the quote-prefix style COM applies so the cell
    # keeps plain default formatting, matching the source rows above it.
    $ws.Cells.Item($r.Row, 3).Value = "'" + $r.C
    $ws.Cells.Item($r.Row, 3).ClearFormats()
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    if ($r.M) {
        $ws.Cells.Item($r.Row, 13).Value = $r.M
    }
}
